$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two username values per the diff
$ws.Range("A2").Value = "adityatest3"
$ws.Range("A4").Value = "adityatest1"

# Move the active selection to B8 (mirrors the sheetView selection change in the diff)
$ws.Range("B8").Select()
